# Slide 2, "Subtitle 2" placeholder: append a new bullet paragraph after
# the existing "...One asmtx for one assessment." line.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# InsertAfter with a leading carriage return starts a new paragraph that
# inherits the preceding paragraph's formatting (bullet, indents, font).
$tr.InsertAfter("`rOne job could include multiple assessments")
